$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sayfa1")

# Duplicate the "Projenin Jenkins'e Entegrayonu" task row, assigning it to
# "Nursema GÜLMEZ" as a new row appended at the bottom of the table (row 28).
$ws.Range("A28").Value = "Projenin Jenkins'e Entegrayonu"
$ws.Range("B28").Value = "Nursema GÜLMEZ"

# Match the borders/formatting used elsewhere in the sheet for this kind of
# row: column A picks up the thin-bottom-border look used by A5:A13/A10,
# column B picks up the thin-box-border look used by B20:B27.
$ws.Range("A10").Copy()
$ws.Range("A28").PasteSpecial(-4122) | Out-Null

$ws.Range("B27").Copy()
$ws.Range("B28").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

$ws.Range("B30").Select() | Out-Null
